$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full desired state of the roster table (rows 2-19), columns A (Name), B (Position), C (Team)
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Mo Bamba", "C", "LA Clippers"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
